$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Column widths: merge E (5) and F (6) into the same custom width
# ------------------------------------------------------------------
$ws.Columns(6).ColumnWidth = $ws.Columns(5).ColumnWidth

# ------------------------------------------------------------------
# Row 6 : E6 number format #,##0 -> #,##0.00 (still empty)
# ------------------------------------------------------------------
$ws.Range("E6").NumberFormat = "#,##0.00"

# ------------------------------------------------------------------
# Row 7 : new regression numbers for the Adult/SCR (patch) row
# ------------------------------------------------------------------
$ws.Range("E7").Value = -1.8039890000000001
$ws.Range("G7").Value = -5.0566409999999999
$ws.Range("F7").NumberFormat = "#,##0.000"
$ws.Range("I7").ClearContents()

# ------------------------------------------------------------------
# Row 8 : blank helper row - fix fonts (Aptos Narrow -> Times New Roman)
# ------------------------------------------------------------------
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("F8").NumberFormat = "#,##0.000"

# ------------------------------------------------------------------
# Row 9 : "Nymph" section header row - add E9/F9 helper cells
# ------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").NumberFormat = "#,##0.00"

$ws.Range("B1").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").NumberFormat = "#,##0.000"

# ------------------------------------------------------------------
# Row 10 : new regression numbers for the Nymph/SCR (patch) row
# ------------------------------------------------------------------
$ws.Range("E10").Value = -0.93405800000000005
$ws.Range("F10").NumberFormat = "#,##0.000"
$ws.Range("I10").ClearContents()

# ------------------------------------------------------------------
# Row 11 : brand new "Clipped PC Analysis" blank row
# ------------------------------------------------------------------
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").ClearContents()

$ws.Range("G10").Copy()
$ws.Range("E11:G11").PasteSpecial(-4122)

# ------------------------------------------------------------------
# View state : window position, selection
# ------------------------------------------------------------------
$wb.Windows(1).Left = 38400
$ws.Range("C5:G10").Select()

Write-Output "edit complete"
